$d = $word.ActiveDocument

function Find-ParagraphByText($doc, $text) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text.TrimEnd() -eq $text) {
            return $p
        }
    }
    return $null
}

# 1) The (until now empty) paragraph that follows
#    "    int getNumLineas() const;" becomes the new
#    "    Linea* getCabezaLineas() const;" declaration, and a fresh
#    blank paragraph is inserted right after it (preserving the blank
#    line that used to separate the getters from the mutators).
$getNumLineasPara = Find-ParagraphByText $d "    int getNumLineas() const;"
$blankAfterGetters = $getNumLineasPara.Next()
$blankAfterGetters.Range.Text = "    Linea* getCabezaLineas() const;"
$blankAfterGetters.Range.InsertParagraphAfter()

# 2) A new blank paragraph is inserted right after "#endif".
$endifPara = Find-ParagraphByText $d "#endif"
$endifPara.Range.InsertParagraphAfter()
